$d = $word.ActiveDocument

# Locate the first paragraph that actually contains text (the leading
# paragraphs in this document are empty placeholder paragraphs that need to
# be removed). A paragraph's Range.Text always includes a trailing
# paragraph-mark character, so a genuinely empty paragraph reports a text
# length of 1.
$firstContentIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.Length -gt 1) {
        $firstContentIndex = $i
        break
    }
}
if ($firstContentIndex -eq 0) {
    $firstContentIndex = 1
}

# Re-anchor the "_GoBack" bookmark to the very start of that first content
# paragraph before removing the preceding empty paragraphs, so that after
# they are deleted the bookmark ends up collapsed at the start of what
# becomes the new first paragraph of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$bmRange = $d.Paragraphs.Item($firstContentIndex).Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the empty paragraphs that precede the first content paragraph.
for ($i = 1; $i -lt $firstContentIndex; $i++) {
    $d.Paragraphs.Item(1).Range.Delete()
}
